$d = $word.ActiveDocument

# --- "Programa" paragraph: split single run into 7 pieces separated by manual line breaks ---
$d.Content.Find.Execute("de sólidos.2)Catálise", $false, $false, $false, $false, $false, $true, 1, $false, "de sólidos.^l2)Catálise", 2)
$d.Content.Find.Execute("de catálise.3)Tensão", $false, $false, $false, $false, $false, $true, 1, $false, "de catálise.^l3)Tensão", 2)
$d.Content.Find.Execute("capilar. Aplicações.4)Classificação", $false, $false, $false, $false, $false, $true, 1, $false, "capilar. Aplicações.^l4)Classificação", 2)
$d.Content.Find.Execute("dispersões coloidais. 5)Interações", $false, $false, $false, $false, $false, $true, 1, $false, "dispersões coloidais. ^l5)Interações", 2)
$d.Content.Find.Execute("em macromoléculas. 6)Estado", $false, $false, $false, $false, $false, $true, 1, $false, "em macromoléculas. ^l6)Estado", 2)
$d.Content.Find.Execute("Coagulação. Aplicações.7)Termodinâmica", $false, $false, $false, $false, $false, $true, 1, $false, "Coagulação. Aplicações.^l7)Termodinâmica", 2)

# --- "Bibliografia" paragraph: split single run into 8 pieces separated by manual line breaks ---
$d.Content.Find.Execute("New York, 19982) BIRDI", $false, $false, $false, $false, $false, $true, 1, $false, "New York, 1998^l2) BIRDI", 2)
$d.Content.Find.Execute("New York, 1997.3) OSHIMA", $false, $false, $false, $false, $false, $true, 1, $false, "New York, 1997.^l3) OSHIMA", 2)
$d.Content.Find.Execute("Oxford, 2006.4) JACOB", $false, $false, $false, $false, $false, $true, 1, $false, "Oxford, 2006.^l4) JACOB", 2)
$d.Content.Find.Execute("Academic, 2010.5) ADAMIAN", $false, $false, $false, $false, $false, $true, 1, $false, "Academic, 2010.^l5) ADAMIAN", 2)
$d.Content.Find.Execute("Materiais, 2002. 6) ADAMSON", $false, $false, $false, $false, $false, $true, 1, $false, "Materiais, 2002. ^l6) ADAMSON", 2)
$d.Content.Find.Execute("John Wiley, 1990.7) SHAW", $false, $false, $false, $false, $false, $true, 1, $false, "John Wiley, 1990.^l7) SHAW", 2)
$d.Content.Find.Execute("185 pp.8)  REGALBUTO", $false, $false, $false, $false, $false, $true, 1, $false, "185 pp.^l8)  REGALBUTO", 2)
